# TC04_Canine_Filter_PrimDisSite-LymphNode.xlsx
# Jenkins-automation pass: drop the `File Type` and `Breed` columns from the
# Neo4j "FilesTab" query stored in B4 of the startup sheet, and restore the
# selection/scroll state to B4 (row 4) as left by the author in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# --- Update the FilesTab query text in B4 -----------------------------
# Original query returned File Name, File Type, Association, Description,
# Format, Size, Case ID, Breed, Diagnosis, Study Code.
# New query drops the `File Type` and `Breed` columns entirely.
$newQuery = "`r`nMATCH (f:file)-->(parent)`r`nWITH DISTINCT f, parent`r`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`r`n MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`r`nWHERE diag.primary_disease_site IN ['Lymph Node']`r`nWITH DISTINCT f, parent, c, demo, diag, s`r`nRETURN coalesce(f.file_name, '') AS ``File Name``, `r`n        coalesce(labels(parent)[0], '') AS ``Association``,`r`n        coalesce(f.file_description, '') AS ``Description``,`r`n        coalesce(f.file_format, '') AS ``Format``,`r`n        coalesce(f.file_size, '') AS ``Size``,`r`n        coalesce(c.case_id, '') AS ``Case ID``, `r`n        coalesce(diag.disease_term,'') AS Diagnosis , `r`n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newQuery

# --- Restore view/selection state --------------------------------------
# The author's session had scrolled the sheet so row 4 is at the top and
# B4 is the active/selected cell (was C2 previously).
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$null = $ws.Range("B4").Select()
